{"js": "// Apply the wording changes to the introductory paragraph:\n//   \"En esta secci\u00f3n vamos a optimizar la estructura de nuestra aplicaci\u00f3n. Por otro\n//   lado, se va a modificar el componente ...\"\n// becomes\n//   \"En esta secci\u00f3n se optimizar\u00e1 la estructura de la aplicaci\u00f3n web. Por otro lado,\n//   se modificar\u00e1 el componente ...\"\n\nconst body = context.document.body;\n\n// 1) \"vamos a optimizar\" -> \"se optimizar\u00e1\"\nconst hit1 = body.search(\"vamos a optimizar\", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\nif (hit1.items.length > 0) {\n  hit1.items[0].insertText(\"se optimizar\u00e1\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"nuestra aplicaci\u00f3n.\" -> \"la aplicaci\u00f3n web.\"\nconst hit2 = body.search(\"nuestra aplicaci\u00f3n.\", { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\nif (hit2.items.length > 0) {\n  hit2.items[0].insertText(\"la aplicaci\u00f3n web.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) \"se va a modificar\" -> \"se modificar\u00e1\"\nconst hit3 = body.search(\"se va a modificar\", { matchCase: true });\nhit3.load(\"items\");\nawait context.sync();\nif (hit3.items.length > 0) {\n  hit3.items[0].insertText(\"se modificar\u00e1\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the wording changes to the introductory paragraph:\n#   \"En esta secci\u00f3n vamos a optimizar la estructura de nuestra aplicaci\u00f3n. Por otro\n#   lado, se va a modificar el componente ...\"\n# becomes\n#   \"En esta secci\u00f3n se optimizar\u00e1 la estructura de la aplicaci\u00f3n web. Por otro lado,\n#   se modificar\u00e1 el componente ...\"\n\n$d = $word.ActiveDocument\n\n# 1) \"vamos a optimizar\" -> \"se optimizar\u00e1\"\n$rng1 = $d.Content\n$rng1.Find.Execute(\"vamos a optimizar\", $false, $false, $false, $false, $false, $true, 1, $false, \"se optimizar\u00e1\", 2)\n\n# 2) \"nuestra aplicaci\u00f3n.\" -> \"la aplicaci\u00f3n web.\"\n$rng2 = $d.Content\n$rng2.Find.Execute(\"nuestra aplicaci\u00f3n.\", $false, $false, $false, $false, $false, $true, 1, $false, \"la aplicaci\u00f3n web.\", 2)\n\n# 3) \"se va a modificar\" -> \"se modificar\u00e1\"\n$rng3 = $d.Content\n$rng3.Find.Execute(\"se va a modificar\", $false, $false, $false, $false, $false, $true, 1, $false, \"se modificar\u00e1\", 2)\n"}
